$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Fix sorting: update the Avg_Time_ms values for the first two data rows
$ws.Range("D2").Value = 0.59056368
$ws.Range("D3").Value = 1.2286066
